# Add a new "release/1.0.0" row to the meta-sheet, marking it as
# deployed/present ("X") for the sit, uat and prod columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "release/1.0.0"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"
